# Insert a new row at row 33 - shifts existing rows 33:90 down to 34:91
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(33).Insert()

# Populate the new row 33 with the new weekly record.
# Constant-across-rows columns (A, B, C, E, F, G, I, Q, R) are filled the
# same as every other data row in this sheet.
$ws.Range("A33").Value = 9
$ws.Range("B33").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 44540
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = 100112022
$ws.Range("G33").Value = "Arveja Verde"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 52
$ws.Range("K33").Value = 15000
$ws.Range("L33").Value = 17000
$ws.Range("M33").Value = 16000
$ws.Range("N33").Value = "`$/saco 25 kilos"
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 640
$ws.Range("Q33").Value = 25
$ws.Range("R33").Value = "Hortaliza"
